# Weekly refresh of the Vega Modelo de Temuco - Coco sheet:
# a new week's record is inserted at the top of the data block (row 15),
# pushing every existing record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current first data row (row 15),
# shifting all existing records (old row 15 .. old row 56) down to
# (new row 16 .. new row 57).
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with this week's record.
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = 'Vega Modelo de Temuco'
$ws.Range("C15").Value = 'La Araucanía'
$ws.Range("D15").Value = 44659
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 'Fruta'
$ws.Range("G15").Value = 100108
$ws.Range("H15").Value = 'Tropicales y subtropicales'
$ws.Range("I15").Value = 100108007
$ws.Range("J15").Value = 'Coco'
$ws.Range("K15").Value = 'Sin especificar'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = 28000
$ws.Range("O15").Value = 28000
$ws.Range("P15").Value = 28000
$ws.Range("Q15").Value = '$/malla 20 unidades'
$ws.Range("R15").Value = 'Perú'
$ws.Range("S15").Value = 1400
$ws.Range("T15").Value = 20
